# Generate Report for Handback
# A new handback round has completed for the "c91ca27e-..." file: its
# handoff/handback timestamps are refreshed on the per-language sheets,
# and the Overview sheet's "Latest HO Xliff Generate Date" is refreshed
# to the newest handoff time across languages (de-de, the last one run).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-20 04:50:55"
$wsZhCn.Range("K3").Value = "2016-08-20 04:51:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-20 04:51:08"
$wsDeDe.Range("K3").Value = "2016-08-20 04:51:31"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-20 04:51:08"
